$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update applicant #2 (row 2) data
$ws.Range("A2").Value = "KST/2022/0067"
$ws.Range("B2").Value = "Gloria"
$ws.Range("D2").Value = "Ansah"

# Update applicant #3 (row 3) data
$ws.Range("A3").Value = "KST/2020/053"
$ws.Range("B3").Value = "Mustapha"
$ws.Range("C3").Value = "Mummin"
$ws.Range("J3").Value = "Ghana"

# Remove the "Institution Code" (R), "Programme Applied Code" (T) and
# "Programme Offered Code" (U) columns entirely
$ws.Columns("R").Delete()
$ws.Columns("S").Delete()
$ws.Columns("S").Delete()

# Update the active selection to C3
$ws.Range("C3").Select()
